$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-coerced to a number
# (trailing-zero / whole-number-looking strings) need an explicit Text format
# applied first so Excel stores them verbatim as strings, matching the source data.
$textForceCells = "D16","D26","D41","D42","D46","D50","D51"
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.876.87'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.904.34'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").Value = '313.13'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Value = '0.5017'
$ws.Range("E7").Value = '  +3.96%  '
$ws.Range("D8").Value = '0.3822'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.07292'
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").Value = '0.9111'
$ws.Range("D11").Value = '20.88'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '0.07639'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '1.915.07'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Value = '5.493'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '6.621'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '91.40'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '0.000008716'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").Value = '27.915.63'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D22").Value = '5.132'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = '154.59'
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").Value = '1.869'
$ws.Range("E25").Value = '  -2.72%  '
$ws.Range("D26").Value = '2.230'
$ws.Range("E26").Value = '  +5.98%  '
$ws.Range("D27").Value = '18.39'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '115.36'
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").Value = '4.937'
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = '0.08989'
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("D31").Value = '3.202'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").Value = '1.239'
$ws.Range("D33").Value = '0.7736'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").Value = '4.651'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = '0.02063'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").Value = '2.566'
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").Value = '1.102'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '0.5529'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '3.014'
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.05282'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").Value = '6.980'
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").Value = '8.550'
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("D43").Value = '0.1525'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '111.28'
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").Value = '10.62'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").Value = '0.4800'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = '1.643'
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").Value = '67.51'
$ws.Range("E49").Value = '  -1.52%  '
$ws.Range("D50").Value = '0.06080'
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").Value = '0.9000'
$ws.Range("E51").Value = '  -0.74%  '
